$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Refresh the timestamp on the previous "batch" of checks (rows 338-351) ---
# These 14 rows were written with a timestamp that gets a tiny precision bump
# (same moment in time, re-serialized).
for ($r = 338; $r -le 351; $r++) {
    $ws.Cells.Item($r, 4).Value2 = 44232.17696618055
}

# --- 2. Append a brand-new "batch" of 14 availability checks as rows 352-365 ---
# Same 14 services/urls cycle used throughout the sheet, all logged at the
# same new run timestamp.
$names = @("Odoo","Blackbox","PowerBI","Dropbox","Odoo","GEE","UtilidadesOdoo","Filtros Dashboard","MapStore","GeoServer","Tomcat","Shiny","Github","EZ Exporter")
$urls = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)

$startRow = 352
$timestamp = 44232.19803154298

for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $names[$i]

    $urlCell = $ws.Cells.Item($r, 2)
    $urlCell.Value = $urls[$i]

    $hashIndex = $urls[$i].IndexOf("#")
    if ($hashIndex -ge 0) {
        $address = $urls[$i].Substring(0, $hashIndex)
        $subAddress = $urls[$i].Substring($hashIndex + 1)
        $ws.Hyperlinks.Add($urlCell, $address, $subAddress) | Out-Null
    } else {
        $ws.Hyperlinks.Add($urlCell, $urls[$i]) | Out-Null
    }
    $urlCell.Style = "Hyperlink"

    $ws.Cells.Item($r, 3).Value = "Disponible"

    $dateCell = $ws.Cells.Item($r, 4)
    $dateCell.Value2 = $timestamp
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
